# Updated SAF displaced CI values
# Applies the edits described in the commit:
#  - C2U CI sheet: fill in Nominal-year "Use Carbon Intensity" values for
#    Sustainable Aviation Fuel (row 14), and fill the High/Low rows (15/16)
#    with formulas copied down from the row above (matching the existing
#    Nominal/High/Low pattern used elsewhere on the sheet).
#  - Extend that sheet's AutoFilter range to include the newly-populated rows.
#  - F2C CI sheet: turn on an AutoFilter criterion that shows only the
#    "Sustainable Aviation Fuel" commodity rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. C2U CI sheet - fill in the SAF "displaced" (conventional jet fuel)
#    carbon-intensity values across the 2025-2045 year columns (G:AA).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("C2U CI")

$values = @(86.64, 85.38, 84.13, 82.87, 81.62, 80.36, 80.36, 80.36, 80.36, 80.36, 80.36, 80.36, 80.36, 80.36, 80.36, 80.36, 80.36, 80.36, 80.36, 80.36, 80.36)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 7 + $i   # column G = 7
    $ws.Cells.Item(14, $col).Value = $values[$i]
}

# Row 15 ("High") = row 14 copied down; row 16 ("Low") = row 15 copied down,
# matching the formula pattern already used by the other commodity blocks
# (column G entered on its own, then H:AA filled right/down together).
$ws.Range("G15").Formula = "=G14"
$ws.Range("H15:AA16").Formula = "=H14"
$ws.Range("G16").Formula = "=G15"

# The sheet's AutoFilter range needs to grow to cover the newly-used rows.
$ws.Range("A1:AA16").AutoFilter()

# ---------------------------------------------------------------------
# 2. F2C CI sheet - filter the Commodity column (B) down to just
#    "Sustainable Aviation Fuel".
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("F2C CI")
$ws2.Range("A1:I61").AutoFilter(2, "Sustainable Aviation Fuel")

# ---------------------------------------------------------------------
# 3. Selection / active-cell bookkeeping to mirror the saved view state.
# ---------------------------------------------------------------------
$wsConversion = $wb.Worksheets.Item("Conversion")
$wsConversion.Range("D16").Select()

$wsF2CCI = $wb.Worksheets.Item("F2C CI")
$wsF2CCI.Range("D72").Select()

$wsC2UAdj = $wb.Worksheets.Item("C2U UO Adjustment")
$wsC2UAdj.Range("C23").Select()

$wsC2UCI = $wb.Worksheets.Item("C2U CI")
$wsC2UCI.Range("F27").Select()
$wsC2UCI.Activate()
